$d = $word.ActiveDocument

$replacements = @(
    @("2024-11-02 Saturday", "2024-11-03 Sunday"),
    @("721×2=1442", "473×3=1419"),
    @("680×7=4760", "781×3=2343"),
    @("568×9=5112", "768×7=5376"),
    @("554×8=4432", "999×5=4995"),
    @("185×3=555", "170×6=1020"),
    @("995×7=6965", "441×8=3528"),
    @("850×8=6800", "822×9=7398"),
    @("941×7=6587", "363×7=2541"),
    @("567×2=1134", "517×9=4653"),
    @("709×9=6381", "561×3=1683"),
    @("532×8=4256", "313×8=2504"),
    @("724×7=5068", "134×2=268"),
    @("664×6=3984", "781×7=5467"),
    @("808×8=6464", "680×8=5440"),
    @("110×7=770", "829×4=3316"),
    @("853×4=3412", "225×5=1125"),
    @("492×5=2460", "857×9=7713"),
    @("461×9=4149", "798×9=7182"),
    @("222×8=1776", "271×5=1355"),
    @("118×4=472", "671×5=3355"),
    @("687×6=4122", "268×8=2144"),
    @("376×3=1128", "653×4=2612"),
    @("448×2=896", "987×7=6909"),
    @("150×4=600", "929×8=7432"),
    @("486×8=3888", "323×8=2584")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
